# "changes for parsing amount efficiently"
# Adds new expense rows (6-14) below the existing data and fixes the
# date format on A5 (was date-only, should be date+time like the other
# data rows).

function Set-TextValue($rng, $val) {
    # Force the cell to be stored as literal text even when the value
    # looks like a number/date/currency amount (e.g. "0", "$1,200",
    # "April 2025"). Restore the default "Normal" style afterwards so we
    # don't leave a stray text-number-format style on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: change the date cell's number format from date-only to
# date+time (style moves from s=3 to s=2).
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 6
$ws.Range("A6").Value = -27277
$ws.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B6").Value = "food"
$ws.Range("C6").Value = 200
Set-TextValue $ws.Range("D6") "April 1825"
$ws.Range("E6").Value = "14:47:48"
$ws.Range("F6").Value = "soya chap at nearby restaurant"

# Row 7
$ws.Range("A7").Value = 45772
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B7").Value = "food"
$ws.Range("C7").Value = 40
Set-TextValue $ws.Range("D7") "April 2025"
$ws.Range("E7").Value = "14:50:02"
$ws.Range("F7").Value = "cold drink"

# Row 8
$ws.Range("A8").Value = 46011
$ws.Range("A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B8").Value = "Gym"
Set-TextValue $ws.Range("C8") '$1,200'
Set-TextValue $ws.Range("D8") "December 2025"
$ws.Range("E8").Value = "14:52:28"
$ws.Range("F8").Value = "gym fees"

# Row 9
$ws.Range("A9").Value = 45772
$ws.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B9").Value = "food"
Set-TextValue $ws.Range("C9") "0"
Set-TextValue $ws.Range("D9") "April 2025"
$ws.Range("E9").Value = "14:54:24"
$ws.Range("F9").Value = "tea"

# Row 10
$ws.Range("A10").Value = 45772
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B10").Value = "tea"
Set-TextValue $ws.Range("C10") '$20,000,000'
Set-TextValue $ws.Range("D10") "April 2025"
$ws.Range("E10").Value = "14:55:12"

# Row 11
$ws.Range("A11").Value = 45772
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B11").Value = "chill"
Set-TextValue $ws.Range("C11") "0"
Set-TextValue $ws.Range("D11") "April 2025"
$ws.Range("E11").Value = "14:55:59"
$ws.Range("F11").Value = "playing games in mall"

# Row 12
$ws.Range("A12").Value = 45772
$ws.Range("A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B12").Value = "food"
Set-TextValue $ws.Range("C12") "0"
Set-TextValue $ws.Range("D12") "April 2025"
$ws.Range("E12").Value = "14:57:12"
$ws.Range("F12").Value = "momos"

# Row 13
$ws.Range("A13").Value = 45144
$ws.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B13").Value = "clothing"
$ws.Range("C13").Value = 75
Set-TextValue $ws.Range("D13") "August 2023"
$ws.Range("E13").Value = "15:00:05"
$ws.Range("F13").Value = "cap"

# Row 14 (date-only format, like the original A5 style)
$ws.Range("A14").Value = 45772
$ws.Range("A14").NumberFormat = "YYYY-MM-DD"
$ws.Range("B14").Value = "food"
$ws.Range("C14").Value = 2500
Set-TextValue $ws.Range("D14") "April 2025"
$ws.Range("E14").Value = "15:01:12"
$ws.Range("F14").Value = "protien powder"
